$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2 ("Tabelle2"): new lookup rows for the google.de examples.
# Shared strings must be created in this exact order (B-col both rows, then
# D-col both rows, then E-col both rows) so they line up with the target
# sharedStrings.xml ordering.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B2").Value = "http://www.google.de/int/a/b"
$ws2.Range("D2").Value = "int/a/b"
$ws2.Range("B3").Value = "http://www.google.de/int/a/b/"
$ws2.Range("D3").Value = "int/a/b/"
$ws2.Range("E2").Value = "www.google.de/int/a/b"
$ws2.Range("E3").Value = "www.google.de/int/a/b/"

$ws2.Hyperlinks.Add($ws2.Range("B2"), "http://www.google.de/int/a/b")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "http://www.google.de/int/a/b/")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "http://www.google.de/int/a/b")
$ws2.Hyperlinks.Add($ws2.Range("E3"), "http://www.google.de/int/a/b/")

$ws2.Range("B2").Style = "Hyperlink"
$ws2.Range("B3").Style = "Hyperlink"
$ws2.Range("E2").Style = "Hyperlink"
$ws2.Range("E3").Style = "Hyperlink"

# Column autosize for the two hyperlink-ish columns that now hold long text.
$ws2.Columns.Item(2).ColumnWidth = 27.833333333333336
$ws2.Columns.Item(5).ColumnWidth = 22.666666666666664

# ---------------------------------------------------------------------------
# Sheet3: renamed to "Requirements" and gets two notes.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Requirements"

$ws3.Range("B3").Value = "einheitliche Überschriften"
$ws3.Range("B4").Value = "Edit / show einheitlich"

# ---------------------------------------------------------------------------
# Selections: Tabelle1 selection moves to B25 (but Tabelle1 is no longer the
# active tab), Tabelle2 selection moves to G9, Requirements becomes the
# active sheet with its selection on B5.
# ---------------------------------------------------------------------------
$ws2.Range("G9").Select()

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B25").Select()

$ws3.Activate()
$ws3.Range("B5").Select()
